# Applies the "Bots" sheet extension: fills in additional bot rows
# (rows 7-21) with "trichterdraws"/"nbl" values using the same cell
# styling (thick left border) already used throughout that block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bots")

function Set-CellNbl($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Value = $val
    # Re-create the "thick left border" styling used by the rest of this
    # block of cells (same visual style as e.g. D10/D11/D12/D13/etc.)
    $cell.Borders(7).Weight = 4
}

Set-CellNbl $ws "C7" "trichterdraws"
Set-CellNbl $ws "C8" "trichterdraws"
Set-CellNbl $ws "C9" "nbl"
Set-CellNbl $ws "F9" "nbl"
Set-CellNbl $ws "B10" "trichterdraws"
Set-CellNbl $ws "C10" "nbl"
Set-CellNbl $ws "F10" "nbl"
Set-CellNbl $ws "A11" "nbl"
Set-CellNbl $ws "B11" "trichterdraws"
Set-CellNbl $ws "C11" "nbl"
Set-CellNbl $ws "E11" "nbl"
Set-CellNbl $ws "F11" "nbl"
Set-CellNbl $ws "A12" "nbl"
Set-CellNbl $ws "B12" "nbl"
Set-CellNbl $ws "C12" "nbl"
Set-CellNbl $ws "E12" "nbl"
Set-CellNbl $ws "F12" "nbl"
Set-CellNbl $ws "A13" "nbl"
Set-CellNbl $ws "B13" "nbl"
Set-CellNbl $ws "E13" "nbl"
Set-CellNbl $ws "F13" "nbl"
Set-CellNbl $ws "A14" "nbl"
Set-CellNbl $ws "B14" "nbl"
Set-CellNbl $ws "D14" "trichterdraws"
Set-CellNbl $ws "E14" "nbl"
Set-CellNbl $ws "F14" "nbl"
Set-CellNbl $ws "A15" "nbl"
Set-CellNbl $ws "B15" "nbl"
Set-CellNbl $ws "D15" "trichterdraws"
Set-CellNbl $ws "E15" "nbl"
Set-CellNbl $ws "A16" "nbl"
Set-CellNbl $ws "B16" "nbl"
Set-CellNbl $ws "D16" "nbl"
Set-CellNbl $ws "E16" "nbl"
Set-CellNbl $ws "A17" "nbl"
Set-CellNbl $ws "B17" "nbl"
Set-CellNbl $ws "D17" "nbl"
Set-CellNbl $ws "E17" "nbl"
Set-CellNbl $ws "A18" "nbl"
Set-CellNbl $ws "D18" "nbl"
Set-CellNbl $ws "E18" "nbl"
Set-CellNbl $ws "E19" "nbl"
Set-CellNbl $ws "E20" "nbl"
Set-CellNbl $ws "E21" "nbl"
